$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-9: update date value from 45184 (2023-09-15) to 45185 (2023-09-16)
$newDate = Get-Date -Year 2023 -Month 9 -Day 16 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
